# Update Leve profit calculations across all job sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) per latest Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1528.5714
$ws.Range("I80").Value = 1066.6666
$ws.Range("J80").Value = 1875
$ws.Range("K80").Value = 3199.9998
$ws.Range("L80").Value = 5625
$ws.Range("M80").Value = -2201.9998
$ws.Range("N80").Value = -7621

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1528.5714
$ws.Range("I83").Value = 1066.6666
$ws.Range("J83").Value = 1875
$ws.Range("K83").Value = 9599.999400000001
$ws.Range("L83").Value = 16875
$ws.Range("M83").Value = -4607.999400000001
$ws.Range("N83").Value = -26859

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3774.8462
$ws.Range("I113").Value = 3703.3333
$ws.Range("J113").Value = 3935.75
$ws.Range("K113").Value = 3703.3333
$ws.Range("L113").Value = 3935.75
$ws.Range("M113").Value = -449.3332999999998
$ws.Range("N113").Value = -10443.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1360.579
$ws.Range("J135").Value = 2011.7142
$ws.Range("L135").Value = 18105.4278
$ws.Range("N135").Value = -23175.4278

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 8334574.5
$ws.Range("I137").Value = 1210.75
$ws.Range("J137").Value = 25001302
$ws.Range("K137").Value = 3632.25
$ws.Range("L137").Value = 75003906
$ws.Range("M137").Value = -1082.25
$ws.Range("N137").Value = -75009006

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 6033.125
$ws.Range("I138").Value = 12066.143
$ws.Range("J138").Value = 2784.577
$ws.Range("K138").Value = 36198.429
$ws.Range("L138").Value = 8353.731
$ws.Range("M138").Value = -31058.429
$ws.Range("N138").Value = -18633.731

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1987.8
$ws.Range("I141").Value = 1513.1666
$ws.Range("K141").Value = 4539.4998
$ws.Range("M141").Value = 640.5002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1309.4
$ws.Range("I2").Value = 986.75
$ws.Range("K2").Value = 986.75
$ws.Range("M2").Value = -873.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 275132.94
$ws.Range("I32").Value = 346231.28
$ws.Range("J32").Value = 17401.5
$ws.Range("K32").Value = 346231.28
$ws.Range("L32").Value = 17401.5
$ws.Range("M32").Value = -345944.28
$ws.Range("N32").Value = -17975.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 65170.47
$ws.Range("I45").Value = 146556.86
$ws.Range("K45").Value = 146556.86
$ws.Range("M45").Value = -146179.86

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 15496.25
$ws.Range("I97").Value = 20629.8
$ws.Range("K97").Value = 20629.8
$ws.Range("M97").Value = -20133.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2417.0908
$ws.Range("I102").Value = 2387.5557
$ws.Range("K102").Value = 2387.5557
$ws.Range("M102").Value = -765.5556999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1309.4
$ws.Range("I116").Value = 986.75
$ws.Range("K116").Value = 986.75
$ws.Range("M116").Value = 1307.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2102.1765
$ws.Range("I132").Value = 1816.0682
$ws.Range("K132").Value = 5448.2046
$ws.Range("M132").Value = -2918.2046

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1309.4
$ws.Range("I3").Value = 986.75
$ws.Range("K3").Value = 986.75
$ws.Range("M3").Value = -872.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2484.8125
$ws.Range("I86").Value = 1487.091
$ws.Range("J86").Value = 4679.8
$ws.Range("K86").Value = 1487.091
$ws.Range("L86").Value = 4679.8
$ws.Range("M86").Value = -364.0909999999999
$ws.Range("N86").Value = -6925.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2484.8125
$ws.Range("I89").Value = 1487.091
$ws.Range("J89").Value = 4679.8
$ws.Range("K89").Value = 7435.455
$ws.Range("L89").Value = 23399
$ws.Range("M89").Value = -1819.455
$ws.Range("N89").Value = -34631

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1522.5769
$ws.Range("I107").Value = 1278.25
$ws.Range("J107").Value = 1913.5
$ws.Range("K107").Value = 1278.25
$ws.Range("L107").Value = 1913.5
$ws.Range("M107").Value = 641.75
$ws.Range("N107").Value = -5753.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 2400.724
$ws.Range("I122").Value = 2474.8262
$ws.Range("K122").Value = 7424.4786
$ws.Range("M122").Value = -4974.4786

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 5670.294
$ws.Range("I55").Value = 5666.6665
$ws.Range("J55").Value = 5671.0713
$ws.Range("K55").Value = 16999.9995
$ws.Range("L55").Value = 17013.2139
$ws.Range("M55").Value = -16822.9995
$ws.Range("N55").Value = -17367.2139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 125000776
$ws.Range("I107").Value = 557
$ws.Range("K107").Value = 1671
$ws.Range("M107").Value = 249

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 1602.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H56").Value = 1000
$ws.Range("I56").Value = 1000
$ws.Range("K56").Value = 1000
$ws.Range("M56").Value = -248

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3174.85
$ws.Range("I22").Value = 1459.1
$ws.Range("J22").Value = 4890.6
$ws.Range("K22").Value = 1459.1
$ws.Range("L22").Value = 4890.6
$ws.Range("M22").Value = -1164.1
$ws.Range("N22").Value = -5480.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 3174.85
$ws.Range("I27").Value = 1459.1
$ws.Range("J27").Value = 4890.6
$ws.Range("K27").Value = 1459.1
$ws.Range("L27").Value = 4890.6
$ws.Range("M27").Value = -1352.1
$ws.Range("N27").Value = -5104.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8036.0835
$ws.Range("J46").Value = 6486
$ws.Range("L46").Value = 6486
$ws.Range("N46").Value = -6862

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2378.0908
$ws.Range("I93").Value = 2219.875
$ws.Range("K93").Value = 2219.875
$ws.Range("M93").Value = -971.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3162.818
$ws.Range("I122").Value = 2825.75
$ws.Range("J122").Value = 3355.4285
$ws.Range("K122").Value = 8477.25
$ws.Range("L122").Value = 10066.2855
$ws.Range("M122").Value = -6027.25
$ws.Range("N122").Value = -14966.2855

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 13320.988
$ws.Range("I136").Value = 47458.184
$ws.Range("K136").Value = 142374.552
$ws.Range("M136").Value = -139824.552

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 450
$ws.Range("I113").Value = 436.2
$ws.Range("K113").Value = 1308.6
$ws.Range("M113").Value = 861.4000000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 50000
$ws.Range("J119").Value = 50000
$ws.Range("L119").Value = 50000
$ws.Range("M119").Value = -59676

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2457.9524
$ws.Range("I122").Value = 1708.1177
$ws.Range("K122").Value = 5124.3531
$ws.Range("M122").Value = -2674.3531

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2606.1875
$ws.Range("I126").Value = 2323.182
$ws.Range("K126").Value = 6969.545999999999
$ws.Range("M126").Value = -4499.545999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 56978.668
$ws.Range("I136").Value = 91434.37
$ws.Range("J136").Value = 2834
$ws.Range("K136").Value = 274303.11
$ws.Range("L136").Value = 8502
$ws.Range("M136").Value = -271753.11
$ws.Range("N136").Value = -13602
